$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2026-02-09 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-02-10 Tuesday", 2)

# Update the multiplication-problem table cells (addressed by row/column so the
# duplicate "15×96=" values can be set independently to their distinct targets)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "63×33="
$t.Cell(1,2).Range.Text  = "37×72="
$t.Cell(1,3).Range.Text  = "39×20="
$t.Cell(1,4).Range.Text  = "73×93="
$t.Cell(1,5).Range.Text  = "47×72="

$t.Cell(5,1).Range.Text  = "81×83="
$t.Cell(5,2).Range.Text  = "98×39="
$t.Cell(5,3).Range.Text  = "87×38="
$t.Cell(5,4).Range.Text  = "37×18="
$t.Cell(5,5).Range.Text  = "14×14="

$t.Cell(10,1).Range.Text = "38×29="
$t.Cell(10,2).Range.Text = "35×42="
$t.Cell(10,3).Range.Text = "68×56="
$t.Cell(10,4).Range.Text = "86×49="
$t.Cell(10,5).Range.Text = "24×61="

$t.Cell(15,1).Range.Text = "42×68="
$t.Cell(15,2).Range.Text = "19×70="
$t.Cell(15,3).Range.Text = "93×76="
$t.Cell(15,4).Range.Text = "12×56="
$t.Cell(15,5).Range.Text = "39×77="

$t.Cell(20,1).Range.Text = "33×19="
$t.Cell(20,2).Range.Text = "29×26="
$t.Cell(20,3).Range.Text = "39×96="
$t.Cell(20,4).Range.Text = "55×70="
$t.Cell(20,5).Range.Text = "76×72="
